$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 116.17647
$ws.Range("I33").Value = 111.333336
$ws.Range("K33").Value = 111.333336
$ws.Range("M33").Value = 117.666664
$ws.Range("H100").Value = 11150131
$ws.Range("I100").Value = 12864998
$ws.Range("K100").Value = 12864998
$ws.Range("M100").Value = -12864457
$ws.Range("H107").Value = 370812.97
$ws.Range("I107").Value = 556000.2
$ws.Range("J107").Value = 438.5
$ws.Range("K107").Value = 556000.2
$ws.Range("L107").Value = 438.5
$ws.Range("M107").Value = -554080.2
$ws.Range("N107").Value = -4278.5
$ws.Range("H129").Value = 1280.68
$ws.Range("J129").Value = 1348.826
$ws.Range("L129").Value = 4046.478
$ws.Range("N129").Value = -14046.478
$ws.Range("H138").Value = 1747.97
$ws.Range("I138").Value = 1056.4524
$ws.Range("J138").Value = 2248.724
$ws.Range("K138").Value = 3169.357199999999
$ws.Range("L138").Value = 6746.172
$ws.Range("M138").Value = 1970.642800000001
$ws.Range("N138").Value = -17026.172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1831.4182
$ws.Range("I61").Value = 1444.8667
$ws.Range("J61").Value = 3570.9
$ws.Range("K61").Value = 1444.8667
$ws.Range("L61").Value = 3570.9
$ws.Range("M61").Value = -1232.8667
$ws.Range("N61").Value = -3994.9
$ws.Range("H122").Value = 13280.333
$ws.Range("H136").Value = 1831.4182
$ws.Range("I136").Value = 1444.8667
$ws.Range("J136").Value = 3570.9
$ws.Range("K136").Value = 4334.6001
$ws.Range("L136").Value = 10712.7
$ws.Range("M136").Value = -1784.6001
$ws.Range("N136").Value = -15812.7
$ws.Range("H139").Value = 40098.69
$ws.Range("J139").Value = 41843.91
$ws.Range("L139").Value = 41843.91
$ws.Range("N139").Value = -52123.91

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2042.4667
$ws.Range("I94").Value = 2092.7
$ws.Range("J94").Value = 1942
$ws.Range("K94").Value = 2092.7
$ws.Range("L94").Value = 1942
$ws.Range("M94").Value = -1641.7
$ws.Range("N94").Value = -2844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2087.7693
$ws.Range("I31").Value = 1147.963
$ws.Range("J31").Value = 4202.3335
$ws.Range("K31").Value = 1147.963
$ws.Range("L31").Value = 4202.3335
$ws.Range("M31").Value = -852.963
$ws.Range("N31").Value = -4792.3335
$ws.Range("H34").Value = 2087.7693
$ws.Range("I34").Value = 1147.963
$ws.Range("J34").Value = 4202.3335
$ws.Range("K34").Value = 1147.963
$ws.Range("L34").Value = 4202.3335
$ws.Range("M34").Value = -945.963
$ws.Range("N34").Value = -4606.3335
$ws.Range("H58").Value = 1401.8914
$ws.Range("I58").Value = 946.4722
$ws.Range("J58").Value = 3041.4
$ws.Range("K58").Value = 946.4722
$ws.Range("L58").Value = 3041.4
$ws.Range("M58").Value = -743.4722
$ws.Range("N58").Value = -3447.4
$ws.Range("H122").Value = 1613
$ws.Range("I122").Value = 1612.5
$ws.Range("J122").Value = 1614
$ws.Range("K122").Value = 4837.5
$ws.Range("L122").Value = 4842
$ws.Range("M122").Value = -2387.5
$ws.Range("N122").Value = -9742
$ws.Range("H132").Value = 2184.45
$ws.Range("I132").Value = 1835.6666
$ws.Range("J132").Value = 3230.8
$ws.Range("K132").Value = 5506.9998
$ws.Range("L132").Value = 9692.400000000001
$ws.Range("M132").Value = -2976.9998
$ws.Range("N132").Value = -14752.4
$ws.Range("H136").Value = 1401.8914
$ws.Range("I136").Value = 946.4722
$ws.Range("J136").Value = 3041.4
$ws.Range("K136").Value = 2839.4166
$ws.Range("L136").Value = 9124.200000000001
$ws.Range("M136").Value = -289.4166
$ws.Range("N136").Value = -14224.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1110.2972
$ws.Range("I5").Value = 667.3461
$ws.Range("J5").Value = 2157.2727
$ws.Range("K5").Value = 2002.0383
$ws.Range("L5").Value = 6471.8181
$ws.Range("M5").Value = -1890.0383
$ws.Range("N5").Value = -6695.8181
$ws.Range("H131").Value = 6668409.5
$ws.Range("J131").Value = 7753885
$ws.Range("L131").Value = 23261655
$ws.Range("N131").Value = -23271735
$ws.Range("H135").Value = 1110.2972
$ws.Range("I135").Value = 667.3461
$ws.Range("J135").Value = 2157.2727
$ws.Range("K135").Value = 6006.1149
$ws.Range("L135").Value = 19415.4543
$ws.Range("M135").Value = -3471.1149
$ws.Range("N135").Value = -24485.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2532.4
$ws.Range("I102").Value = 2854
$ws.Range("J102").Value = 2050
$ws.Range("K102").Value = 2854
$ws.Range("L102").Value = 2050
$ws.Range("M102").Value = -1232
$ws.Range("N102").Value = -5294
$ws.Range("H122").Value = 696340
$ws.Range("I122").Value = 1390490.5
$ws.Range("J122").Value = 2189.5
$ws.Range("K122").Value = 4171471.5
$ws.Range("L122").Value = 6568.5
$ws.Range("M122").Value = -4169021.5
$ws.Range("N122").Value = -11468.5
$ws.Range("H123").Value = 9807.75
$ws.Range("J123").Value = 9807.75
$ws.Range("L123").Value = 9807.75
$ws.Range("N123").Value = -14707.75
$ws.Range("H126").Value = 1841.6279
$ws.Range("I126").Value = 1465.3889
$ws.Range("J126").Value = 2112.52
$ws.Range("K126").Value = 4396.1667
$ws.Range("L126").Value = 6337.559999999999
$ws.Range("M126").Value = -1926.1667
$ws.Range("N126").Value = -11277.56
$ws.Range("H132").Value = 4322.3335
$ws.Range("I132").Value = 4388.2104
$ws.Range("J132").Value = 4165.875
$ws.Range("K132").Value = 13164.6312
$ws.Range("L132").Value = 12497.625
$ws.Range("M132").Value = -10634.6312
$ws.Range("N132").Value = -17557.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3062.75
$ws.Range("I40").Value = 2683.1667
$ws.Range("J40").Value = 3290.5
$ws.Range("K40").Value = 2683.1667
$ws.Range("L40").Value = 3290.5
$ws.Range("M40").Value = -2547.1667
$ws.Range("N40").Value = -3562.5
$ws.Range("H100").Value = 3115.9473
$ws.Range("I100").Value = 2740.6
$ws.Range("J100").Value = 3250
$ws.Range("K100").Value = 2740.6
$ws.Range("L100").Value = 3250
$ws.Range("M100").Value = -2199.6
$ws.Range("N100").Value = -4332
$ws.Range("H132").Value = 5328.6978
$ws.Range("I132").Value = 4944.9653
$ws.Range("J132").Value = 6123.5713
$ws.Range("K132").Value = 14834.8959
$ws.Range("L132").Value = 18370.7139
$ws.Range("M132").Value = -12304.8959
$ws.Range("N132").Value = -23430.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 27604
$ws.Range("J121").Value = 27604
$ws.Range("L121").Value = 27604
$ws.Range("N121").Value = -31098
$ws.Range("H123").Value = 33749.562
$ws.Range("J123").Value = 33749.562
$ws.Range("L123").Value = 33749.562
$ws.Range("N123").Value = -43549.562
